$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 (shifts existing rows 16.. down by one,
# copying formatting such as the date number format from the row above).
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new record.
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(16, 3).Value = "La Araucanía"
$ws.Cells.Item(16, 4).Value = "10/31/2023"
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = 300000000
$ws.Cells.Item(16, 7).Value = "Espárragos"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 250
$ws.Cells.Item(16, 11).Value = 1600
$ws.Cells.Item(16, 12).Value = 1600
$ws.Cells.Item(16, 13).Value = 1600
$ws.Cells.Item(16, 14).Value = "`$/kilo"
$ws.Cells.Item(16, 15).Value = "Región del Maule"
$ws.Cells.Item(16, 16).Value = 1600
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(16, 18).Value = "Hortaliza"
